$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestOutput")

# Header row shifts left by one column (the old "test_case" column B is
# removed) and A1's label changes from test_suite_id -> test_case_id.
$ws.Cells.Item(1, 1).Value = "test_case_id"
$ws.Cells.Item(1, 2).Value = "pks"
$ws.Cells.Item(1, 3).Value = "id"
$ws.Cells.Item(1, 4).Value = "name"
$ws.Cells.Item(1, 5).Value = "description"
$ws.Cells.Item(1, 6).Value = "tags"

# Drop the now-unused trailing column G so the used range is A1:F1.
$ws.Cells.Item(1, 7).ClearContents()
